# COVID-19 Bolivia BD_dpto.xlsx update — "Actualizado 9 de agosto de 2020"
# Appends 7 new daily rows (2020-08-03 .. 2020-08-09, serials 44046-44052)
# to each of the 9 department sheets (rows 148-154), matching the existing
# table formatting by copying row 146 (format only) as a template.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Beni" ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Beni"
$ws.Range("C148").Value = 19
$ws.Range("D148").Value = 2
$ws.Range("E148").Value = 56
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Beni"
$ws.Range("C149").Value = 31
$ws.Range("D149").Value = 1
$ws.Range("E149").Value = 13
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Beni"
$ws.Range("C150").Value = 28
$ws.Range("D150").Value = 3
$ws.Range("E150").Value = 38
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Beni"
$ws.Range("C151").Value = 32
$ws.Range("D151").Value = 2
$ws.Range("E151").Value = 13
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Beni"
$ws.Range("C152").Value = 31
$ws.Range("D152").Value = 2
$ws.Range("E152").Value = 31
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Beni"
$ws.Range("C153").Value = 17
$ws.Range("D153").Value = 2
$ws.Range("E153").Value = 2
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Beni"
$ws.Range("C154").Value = 12
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 4
[void]$ws.Range("E153").Select()

# --- Sheet 2: "Cochabamba" ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Cochabamba"
$ws.Range("C148").Value = 189
$ws.Range("D148").Value = 15
$ws.Range("E148").Value = 54
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Cochabamba"
$ws.Range("C149").Value = 216
$ws.Range("D149").Value = 17
$ws.Range("E149").Value = 108
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Cochabamba"
$ws.Range("C150").Value = 182
$ws.Range("D150").Value = 15
$ws.Range("E150").Value = 74
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Cochabamba"
$ws.Range("C151").Value = 96
$ws.Range("D151").Value = 14
$ws.Range("E151").Value = 78
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Cochabamba"
$ws.Range("C152").Value = 83
$ws.Range("D152").Value = 12
$ws.Range("E152").Value = 82
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Cochabamba"
$ws.Range("C153").Value = 43
$ws.Range("D153").Value = 10
$ws.Range("E153").Value = 83
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Cochabamba"
$ws.Range("C154").Value = 151
$ws.Range("D154").Value = 5
$ws.Range("E154").Value = 147
[void]$ws.Range("E153").Select()

# --- Sheet 3: "Chuquisaca" ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Chuquisaca"
$ws.Range("C148").Value = 58
$ws.Range("D148").Value = 9
$ws.Range("E148").Value = 8
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Chuquisaca"
$ws.Range("C149").Value = 62
$ws.Range("D149").Value = 4
$ws.Range("E149").Value = 5
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Chuquisaca"
$ws.Range("C150").Value = 76
$ws.Range("D150").Value = 6
$ws.Range("E150").Value = 59
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Chuquisaca"
$ws.Range("C151").Value = 99
$ws.Range("D151").Value = 9
$ws.Range("E151").Value = 58
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Chuquisaca"
$ws.Range("C152").Value = 126
$ws.Range("D152").Value = 13
$ws.Range("E152").Value = 47
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Chuquisaca"
$ws.Range("C153").Value = 99
$ws.Range("D153").Value = 13
$ws.Range("E153").Value = 61
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Chuquisaca"
$ws.Range("C154").Value = 51
$ws.Range("D154").Value = 8
$ws.Range("E154").Value = 63
[void]$ws.Range("E153").Select()

# --- Sheet 4: "La Paz" ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "La Paz"
$ws.Range("C148").Value = 876
$ws.Range("D148").Value = 1
$ws.Range("E148").Value = 0
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "La Paz"
$ws.Range("C149").Value = 416
$ws.Range("D149").Value = 21
$ws.Range("E149").Value = 26
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "La Paz"
$ws.Range("C150").Value = 909
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 79
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "La Paz"
$ws.Range("C151").Value = 551
$ws.Range("D151").Value = 16
$ws.Range("E151").Value = 62
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "La Paz"
$ws.Range("C152").Value = 906
$ws.Range("D152").Value = 6
$ws.Range("E152").Value = 39
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "La Paz"
$ws.Range("C153").Value = 624
$ws.Range("D153").Value = 9
$ws.Range("E153").Value = 63
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "La Paz"
$ws.Range("C154").Value = 501
$ws.Range("D154").Value = 10
$ws.Range("E154").Value = 81
[void]$ws.Range("E153").Select()

# --- Sheet 5: "Oruro" ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Oruro"
$ws.Range("C148").Value = 48
$ws.Range("D148").Value = 4
$ws.Range("E148").Value = 11
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Oruro"
$ws.Range("C149").Value = 149
$ws.Range("D149").Value = 6
$ws.Range("E149").Value = 37
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Oruro"
$ws.Range("C150").Value = 70
$ws.Range("D150").Value = 3
$ws.Range("E150").Value = 9
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Oruro"
$ws.Range("C151").Value = 67
$ws.Range("D151").Value = 4
$ws.Range("E151").Value = 20
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Oruro"
$ws.Range("C152").Value = 12
$ws.Range("D152").Value = 0
$ws.Range("E152").Value = 5
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Oruro"
$ws.Range("C153").Value = 78
$ws.Range("D153").Value = 1
$ws.Range("E153").Value = 20
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Oruro"
$ws.Range("C154").Value = 21
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 2
[void]$ws.Range("D153").Select()

# --- Sheet 6: "Pando" ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Pando"
$ws.Range("C148").Value = 3
$ws.Range("D148").Value = 6
$ws.Range("E148").Value = 32
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Pando"
$ws.Range("C149").Value = 28
$ws.Range("D149").Value = 0
$ws.Range("E149").Value = 0
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Pando"
$ws.Range("C150").Value = 36
$ws.Range("D150").Value = 5
$ws.Range("E150").Value = 0
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Pando"
$ws.Range("C151").Value = 22
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 0
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Pando"
$ws.Range("C152").Value = 32
$ws.Range("D152").Value = 2
$ws.Range("E152").Value = 1
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Pando"
$ws.Range("C153").Value = 38
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 1
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Pando"
$ws.Range("C154").Value = 35
$ws.Range("D154").Value = 3
$ws.Range("E154").Value = 0
[void]$ws.Range("C153").Select()

# --- Sheet 7: "Potosí" ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Potosí"
$ws.Range("C148").Value = 51
$ws.Range("D148").Value = 8
$ws.Range("E148").Value = 3
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Potosí"
$ws.Range("C149").Value = 43
$ws.Range("D149").Value = 3
$ws.Range("E149").Value = 15
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Potosí"
$ws.Range("C150").Value = 77
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 68
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Potosí"
$ws.Range("C151").Value = 13
$ws.Range("D151").Value = 2
$ws.Range("E151").Value = 0
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Potosí"
$ws.Range("C152").Value = 32
$ws.Range("D152").Value = 0
$ws.Range("E152").Value = 6
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Potosí"
$ws.Range("C153").Value = 28
$ws.Range("D153").Value = 5
$ws.Range("E153").Value = 17
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Potosí"
$ws.Range("C154").Value = 40
$ws.Range("D154").Value = 4
$ws.Range("E154").Value = 6
[void]$ws.Range("F154").Select()

# --- Sheet 8: "Santa Cruz" ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Santa Cruz"
$ws.Range("C148").Value = 324
$ws.Range("D148").Value = 27
$ws.Range("E148").Value = 155
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Santa Cruz"
$ws.Range("C149").Value = 463
$ws.Range("D149").Value = 35
$ws.Range("E149").Value = 615
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Santa Cruz"
$ws.Range("C150").Value = 322
$ws.Range("D150").Value = 31
$ws.Range("E150").Value = 658
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Santa Cruz"
$ws.Range("C151").Value = 378
$ws.Range("D151").Value = 29
$ws.Range("E151").Value = 692
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Santa Cruz"
$ws.Range("C152").Value = 229
$ws.Range("D152").Value = 19
$ws.Range("E152").Value = 536
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Santa Cruz"
$ws.Range("C153").Value = 236
$ws.Range("D153").Value = 18
$ws.Range("E153").Value = 507
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Santa Cruz"
$ws.Range("C154").Value = 133
$ws.Range("D154").Value = 21
$ws.Range("E154").Value = 563
[void]$ws.Range("F154").Select()

# --- Sheet 9: "Tarija" ---
$ws = $wb.Worksheets.Item(9)
$ws.Range("A146:E146").Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Range("A148").Value = 44046
$ws.Range("B148").Value = "Tarija"
$ws.Range("C148").Value = 125
$ws.Range("D148").Value = 3
$ws.Range("E148").Value = 35
$ws.Range("A146:E146").Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Range("A149").Value = 44047
$ws.Range("B149").Value = "Tarija"
$ws.Range("C149").Value = 107
$ws.Range("D149").Value = 5
$ws.Range("E149").Value = 61
$ws.Range("A146:E146").Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Range("A150").Value = 44048
$ws.Range("B150").Value = "Tarija"
$ws.Range("C150").Value = 80
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 62
$ws.Range("A146:E146").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Range("A151").Value = 44049
$ws.Range("B151").Value = "Tarija"
$ws.Range("C151").Value = 24
$ws.Range("D151").Value = 4
$ws.Range("E151").Value = 13
$ws.Range("A146:E146").Copy()
$ws.Range("A152:E152").PasteSpecial(-4122)
$ws.Range("A152").Value = 44050
$ws.Range("B152").Value = "Tarija"
$ws.Range("C152").Value = 17
$ws.Range("D152").Value = 5
$ws.Range("E152").Value = 19
$ws.Range("A146:E146").Copy()
$ws.Range("A153:E153").PasteSpecial(-4122)
$ws.Range("A153").Value = 44051
$ws.Range("B153").Value = "Tarija"
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 5
$ws.Range("E153").Value = 11
$ws.Range("A146:E146").Copy()
$ws.Range("A154:E154").PasteSpecial(-4122)
$ws.Range("A154").Value = 44052
$ws.Range("B154").Value = "Tarija"
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 2
$ws.Range("E154").Value = 38
[void]$ws.Range("C154").Select()

$excel.CutCopyMode = $false
